$d = $word.ActiveDocument

function Replace-Text($range, $find, $replace) {
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 0, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# --- "Programa" paragraph: split the single run into three, joined by manual line breaks ---
Replace-Text $d.Content `
    "Proposição pelos alunos de startup de base tecnológica A proposta de startup é acompanhada por uma equipe de mentores, coordenada pelo professor da disciplina. Os mentores realizam apresentações sobre" `
    "Proposição pelos alunos de startup de base tecnológica ^lA proposta de startup é acompanhada por uma equipe de mentores, coordenada pelo professor da disciplina. ^lOs mentores realizam apresentações sobre"

# --- "Avaliação" paragraph, "Método" run: split into two sentences with a manual line break ---
Replace-Text $d.Content `
    "Atividades docentes: Mentoria, palestras e seminários.Atividades discentes:" `
    "Atividades docentes: Mentoria, palestras e seminários.^lAtividades discentes:"

# --- "Avaliação" paragraph, "Critério" run: split into two sentences with a manual line break ---
Replace-Text $d.Content `
    "qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.Nota de projeto" `
    "qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.^lNota de projeto"

# --- "Bibliografia" paragraph: split the run into six references, separated by blank lines (two manual breaks) ---
Replace-Text $d.Content `
    "Rio de Janeiro: Campus, 2010.INPI." `
    "Rio de Janeiro: Campus, 2010.^l^lINPI."

Replace-Text $d.Content `
    "Consultado em: junho de 2015.KUMAR," `
    "Consultado em: junho de 2015.^l^lKUMAR,"

Replace-Text $d.Content `
    "New Jersey: John Willey and Sons, 2013.MALHOTRA," `
    "New Jersey: John Willey and Sons, 2013.^l^lMALHOTRA,"

Replace-Text $d.Content `
    "Porto Alegre: Bookman, 2006.ROMEIRO FILHO" `
    "Porto Alegre: Bookman, 2006.^l^lROMEIRO FILHO"

Replace-Text $d.Content `
    "Rio de Janeiro: Campus, 2010.ROZENFELD," `
    "Rio de Janeiro: Campus, 2010.^l^lROZENFELD,"

Write-Output "done"
